$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fill in values for rows 5 and 6 (B, C, D columns) matching rows 3/4
$ws.Range("B5").Value = "Чеботарёв Н.Ю."
$ws.Range("D5").Value = "nikolay_chebotaryov"
$ws.Range("B6").Value = "Чеботарёв Н.Ю."
$ws.Range("D6").Value = "nikolay_chebotaryov"

# Add hyperlink covering C5:C6 with display text, which also sets the cell value
$ws.Hyperlinks.Add($ws.Range("C5:C6"), "mailto:chebotaryovnik@gmail.com", "", "", "chebotaryovnik@gmail.com")

# Update the active selection
$ws.Range("C9").Select()
